$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 959.45
$ws.Range("I18").Value = 764.35297
$ws.Range("J18").Value = 2065
$ws.Range("K18").Value = 764.35297
$ws.Range("L18").Value = 2065
$ws.Range("M18").Value = -480.35297
$ws.Range("N18").Value = -2633
$ws.Range("H40").Value = 1980.9524
$ws.Range("I40").Value = 1850.125
$ws.Range("J40").Value = 2061.4614
$ws.Range("K40").Value = 1850.125
$ws.Range("L40").Value = 2061.4614
$ws.Range("M40").Value = -1675.125
$ws.Range("N40").Value = -2411.4614
$ws.Range("H64").Value = 41281.92
$ws.Range("I64").Value = 127512.5
$ws.Range("J64").Value = 2957.2222
$ws.Range("K64").Value = 127512.5
$ws.Range("L64").Value = 2957.2222
$ws.Range("M64").Value = -127264.5
$ws.Range("N64").Value = -3453.2222
$ws.Range("H67").Value = 41281.92
$ws.Range("I67").Value = 127512.5
$ws.Range("J67").Value = 2957.2222
$ws.Range("K67").Value = 127512.5
$ws.Range("L67").Value = 2957.2222
$ws.Range("M67").Value = -126654.5
$ws.Range("N67").Value = -4673.2222
$ws.Range("H103").Value = 1017.8333
$ws.Range("I103").Value = 2787.5
$ws.Range("J103").Value = 133
$ws.Range("K103").Value = 8362.5
$ws.Range("L103").Value = 399
$ws.Range("M103").Value = -7776.5
$ws.Range("N103").Value = -1571
$ws.Range("H113").Value = 168834.17
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -8508

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29790.463
$ws.Range("I32").Value = 7495.915
$ws.Range("J32").Value = 194212.75
$ws.Range("K32").Value = 7495.915
$ws.Range("L32").Value = 194212.75
$ws.Range("M32").Value = -7208.915
$ws.Range("N32").Value = -194786.75
$ws.Range("H45").Value = 76235
$ws.Range("I45").Value = 144684.72
$ws.Range("J45").Value = 7785.2856
$ws.Range("K45").Value = 144684.72
$ws.Range("L45").Value = 7785.2856
$ws.Range("M45").Value = -144307.72
$ws.Range("N45").Value = -8539.285599999999
$ws.Range("H102").Value = 114395.445
$ws.Range("I102").Value = 169946.5
$ws.Range("J102").Value = 3293.3333
$ws.Range("K102").Value = 169946.5
$ws.Range("L102").Value = 3293.3333
$ws.Range("M102").Value = -168324.5
$ws.Range("N102").Value = -6537.3333
$ws.Range("H122").Value = 1288.5358
$ws.Range("I122").Value = 1084.3
$ws.Range("K122").Value = 3252.9
$ws.Range("M122").Value = -802.8999999999996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2200
$ws.Range("I99").Value = 3580
$ws.Range("J99").Value = 2118.8235
$ws.Range("K99").Value = 3580
$ws.Range("L99").Value = 2118.8235
$ws.Range("M99").Value = -2082
$ws.Range("N99").Value = -5114.8235
$ws.Range("H107").Value = 500005000
$ws.Range("I107").Value = 500005000
$ws.Range("K107").Value = 500005000
$ws.Range("M107").Value = -500003080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2425.5386
$ws.Range("I62").Value = 2218.8572
$ws.Range("J62").Value = 2666.6667
$ws.Range("K62").Value = 2218.8572
$ws.Range("L62").Value = 2666.6667
$ws.Range("M62").Value = -1594.8572
$ws.Range("N62").Value = -3914.6667
$ws.Range("H65").Value = 2425.5386
$ws.Range("I65").Value = 2218.8572
$ws.Range("J65").Value = 2666.6667
$ws.Range("K65").Value = 11094.286
$ws.Range("L65").Value = 13333.3335
$ws.Range("M65").Value = -7974.286
$ws.Range("N65").Value = -19573.3335
$ws.Range("H74").Value = 28367.5
$ws.Range("J74").Value = 28367.5
$ws.Range("L74").Value = 28367.5
$ws.Range("N74").Value = -30115.5
$ws.Range("H77").Value = 28367.5
$ws.Range("J77").Value = 28367.5
$ws.Range("L77").Value = 85102.5
$ws.Range("N77").Value = -93838.5
$ws.Range("H99").Value = 2982.3333
$ws.Range("J99").Value = 2782.8
$ws.Range("L99").Value = 2782.8
$ws.Range("N99").Value = -5778.8
$ws.Range("H107").Value = 1222
$ws.Range("I107").Value = 1883.2
$ws.Range("J107").Value = 808.75
$ws.Range("K107").Value = 1883.2
$ws.Range("L107").Value = 808.75
$ws.Range("M107").Value = 36.79999999999995
$ws.Range("N107").Value = -4648.75
$ws.Range("H122").Value = 1002.5
$ws.Range("J122").Value = 1099
$ws.Range("L122").Value = 3297
$ws.Range("N122").Value = -8197
$ws.Range("H126").Value = 2982.3333
$ws.Range("J126").Value = 2782.8
$ws.Range("L126").Value = 8348.400000000001
$ws.Range("N126").Value = -13288.4
$ws.Range("H132").Value = 3920.2273
$ws.Range("I132").Value = 3855.1052
$ws.Range("J132").Value = 4332.6665
$ws.Range("K132").Value = 11565.3156
$ws.Range("L132").Value = 12997.9995
$ws.Range("M132").Value = -9035.3156
$ws.Range("N132").Value = -18057.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 828.76
$ws.Range("J131").Value = 843.0516
$ws.Range("L131").Value = 2529.1548
$ws.Range("N131").Value = -12609.1548

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 160748.31
$ws.Range("I70").Value = 255601
$ws.Range("K70").Value = 255601
$ws.Range("M70").Value = -255331
$ws.Range("H73").Value = 160748.31
$ws.Range("I73").Value = 255601
$ws.Range("K73").Value = 255601
$ws.Range("M73").Value = -254665
$ws.Range("H107").Value = 842164.3
$ws.Range("I107").Value = 330.22223
$ws.Range("J107").Value = 3367666.8
$ws.Range("K107").Value = 330.22223
$ws.Range("L107").Value = 3367666.8
$ws.Range("M107").Value = 1589.77777
$ws.Range("N107").Value = -3371506.8
$ws.Range("H108").Value = 32000
$ws.Range("J108").Value = 32000
$ws.Range("L108").Value = 32000
$ws.Range("N108").Value = -39680
$ws.Range("H113").Value = 1499.3334
$ws.Range("I113").Value = 1681.6666
$ws.Range("J113").Value = 1408.1666
$ws.Range("K113").Value = 1681.6666
$ws.Range("L113").Value = 1408.1666
$ws.Range("M113").Value = 488.3334
$ws.Range("N113").Value = -5748.1666
$ws.Range("H122").Value = 1955.826
$ws.Range("I122").Value = 1229.8572
$ws.Range("J122").Value = 3085.111
$ws.Range("K122").Value = 3689.5716
$ws.Range("L122").Value = 9255.332999999999
$ws.Range("M122").Value = -1239.5716
$ws.Range("N122").Value = -14155.333
$ws.Range("H126").Value = 3711.5
$ws.Range("I126").Value = 3698.8572
$ws.Range("J126").Value = 3800
$ws.Range("K126").Value = 11096.5716
$ws.Range("L126").Value = 11400
$ws.Range("M126").Value = -8626.571599999999
$ws.Range("N126").Value = -16340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5311.9375
$ws.Range("I7").Value = 3725.875
$ws.Range("J7").Value = 6898
$ws.Range("K7").Value = 3725.875
$ws.Range("L7").Value = 6898
$ws.Range("M7").Value = -3613.875
$ws.Range("N7").Value = -7122
$ws.Range("H40").Value = 48486.816
$ws.Range("I40").Value = 103990.6
$ws.Range("K40").Value = 103990.6
$ws.Range("M40").Value = -103854.6
$ws.Range("H46").Value = 633206.25
$ws.Range("I46").Value = 332.22223
$ws.Range("J46").Value = 1446901.4
$ws.Range("K46").Value = 332.22223
$ws.Range("L46").Value = 1446901.4
$ws.Range("M46").Value = -144.22223
$ws.Range("N46").Value = -1447277.4
$ws.Range("H61").Value = 1898.625
$ws.Range("I61").Value = 1888.2222
$ws.Range("J61").Value = 1912
$ws.Range("K61").Value = 1888.2222
$ws.Range("L61").Value = 1912
$ws.Range("M61").Value = -1686.2222
$ws.Range("N61").Value = -2316
$ws.Range("H68").Value = 5746.6665
$ws.Range("I68").Value = 3400
$ws.Range("J68").Value = 6040
$ws.Range("K68").Value = 3400
$ws.Range("L68").Value = 6040
$ws.Range("M68").Value = -2651
$ws.Range("N68").Value = -7538
$ws.Range("H69").Value = 30387.334
$ws.Range("J69").Value = 30387.334
$ws.Range("L69").Value = 30387.334
$ws.Range("N69").Value = -32009.334
$ws.Range("H71").Value = 5746.6665
$ws.Range("I71").Value = 3400
$ws.Range("J71").Value = 6040
$ws.Range("K71").Value = 17000
$ws.Range("L71").Value = 30200
$ws.Range("M71").Value = -13256
$ws.Range("N71").Value = -37688
$ws.Range("H72").Value = 30387.334
$ws.Range("J72").Value = 30387.334
$ws.Range("L72").Value = 91162.00199999999
$ws.Range("N72").Value = -99274.00199999999
$ws.Range("H93").Value = 2208.7727
$ws.Range("J93").Value = 2165.4285
$ws.Range("L93").Value = 2165.4285
$ws.Range("N93").Value = -4661.4285
$ws.Range("H113").Value = 1898.625
$ws.Range("I113").Value = 1888.2222
$ws.Range("J113").Value = 1912
$ws.Range("K113").Value = 1888.2222
$ws.Range("L113").Value = 1912
$ws.Range("M113").Value = 281.7778000000001
$ws.Range("N113").Value = -6252
$ws.Range("H122").Value = 2721.2
$ws.Range("I122").Value = 2721.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8163.599999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5713.599999999999
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 5311.9375
$ws.Range("I126").Value = 3725.875
$ws.Range("J126").Value = 6898
$ws.Range("K126").Value = 11177.625
$ws.Range("L126").Value = 20694
$ws.Range("M126").Value = -8707.625
$ws.Range("N126").Value = -25634
$ws.Range("H132").Value = 5858.8237
$ws.Range("I132").Value = 6250.3
$ws.Range("K132").Value = 18750.9
$ws.Range("M132").Value = -16220.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 651.8461
$ws.Range("I107").Value = 489.33334
$ws.Range("J107").Value = 1017.5
$ws.Range("K107").Value = 1468.00002
$ws.Range("L107").Value = 3052.5
$ws.Range("M107").Value = 451.9999800000001
$ws.Range("N107").Value = -6892.5
$ws.Range("H122").Value = 3685.625
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 3855
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 11565
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -16465
$ws.Range("H132").Value = 2030.1936
$ws.Range("I132").Value = 2090.0466
$ws.Range("J132").Value = 1894.7368
$ws.Range("K132").Value = 6270.139800000001
$ws.Range("L132").Value = 5684.2104
$ws.Range("M132").Value = -3740.139800000001
$ws.Range("N132").Value = -10744.2104
